# Update countries & provincias Spain
#
# 1. Brasil (row 12) - refreshed daily stats.
# 2. Sri Lanka (row 103) - refreshed daily stats.
# 3. A new "Guatemala" entry is inserted right after "Sri Lanka", with
#    fresh data; the rows that used to follow Sri Lanka
#    (Consejo Danes para los Refugiados, Somalia, Uruguay, Burkina Faso)
#    shift down by one row, and the old "Guatemala" row (which used to sit
#    right before "Sudan") is removed, leaving "Sudan" where it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Brasil (row 12): refresh daily case counts ---
$ws.Range("B12").Value = 97100
$ws.Range("C12").Value = 541
$ws.Range("E12").Value = 49402
$ws.Range("G12").Value = 11
$ws.Range("H12").Value = 6761

# --- Sri Lanka (row 103): refresh daily case counts ---
$ws.Range("B103").Value = 705
$ws.Range("C103").Value = 3
$ws.Range("E103").Value = 526

# --- Insert a single new row right after Sri Lanka (row 103) for
#     Guatemala, pushing the following countries down one row ---
$ws.Rows.Item(104).Insert()

$ws.Range("A104").Value = "Guatemala"
$ws.Range("B104").Value = 688
$ws.Range("C104").Value = 44
$ws.Range("D104").Value = 72
$ws.Range("E104").Value = 599
$ws.Range("F104").Value = 5
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 17

# --- Remove the old Guatemala row, which (after the insert above) is now
#     duplicated at row 109, right before "Sudan" ---
$ws.Rows.Item(109).Delete()
